# Update loading_percent values for rows A=0..23 (sheet rows 2..25), columns B..N
# (column E remains 0, matching the original data)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
        @(14.80612595998577, 5.688914384200675, 7.433330301160994, 0, 41.07355706935063, 48.38296573053749, 19.28609938388824, 30.42795352077218, 10.52040685865705, 12.20453708993182, 11.8460302092484, 16.43345541587372, 22.46450857345616),
        @(14.65312028654114, 5.571236714660224, 7.428254365287355, 0, 41.14281567594985, 48.45394970695931, 19.33232995972031, 30.5078490732858, 10.53636907237344, 12.09539462162928, 11.85727002920242, 16.42075543751854, 22.52615590170458),
        @(14.56138822428877, 5.496677674401181, 7.425949518872491, 0, 41.19284837216055, 48.50812593215856, 19.36337936735459, 30.56150165350656, 10.54672707243143, 12.0302876638906, 11.86551545490169, 16.4152507495921, 22.56581841370652),
        @(14.52460358655552, 5.46573482353373, 7.425215609953161, 0, 41.21512347449497, 48.53286243492127, 19.37670212345568, 30.58452095882555, 10.55108852718973, 12.00426205444331, 11.86921414409781, 16.41358710288995, 22.58243779254768),
        @(14.51853273052283, 5.460563606057653, 7.425106180375438, 0, 41.21893612023024, 48.53713038788936, 19.37895481736408, 30.58841307188818, 10.55182123987783, 11.99997183904174, 11.86984877722163, 16.41334593998296, 22.58522504081321),
        @(14.56088966462267, 5.496262606111887, 7.425938788183983, 0, 41.19314114637421, 48.50844877617851, 19.36355633034744, 30.56180742207055, 10.5467853231957, 12.02993459032735, 11.86556396490223, 16.41522596278715, 22.56604069813334),
        @(14.75293451002885, 5.648829082529201, 7.431412369401966, 0, 41.09587843568688, 48.4052404256262, 19.30148689674156, 30.45454732343796, 10.52579523077006, 12.16652513679476, 11.84962710641561, 16.42860234861351, 22.48538949094934),
        @(15.14528925234744, 5.928906823343101, 7.448532459720206, 0, 40.96477142307958, 48.28704495053251, 19.20090029831428, 30.28069624075923, 10.48903665489589, 12.44830621231038, 11.82901210738993, 16.47289088492826, 22.34154296837627),
        @(15.44068391235287, 6.122021859309949, 7.464927308309999, 0, 40.90485466967257, 48.25170116573216, 19.13987654650829, 30.17523058678488, 10.46468959742871, 12.6621643978373, 11.82031223367692, 16.5162401967933, 22.24450072852202),
        @(15.57608347794855, 6.206936037315407, 7.473196761789046, 0, 40.88550722380215, 48.24682121339943, 19.11491006251327, 30.13208939397441, 10.45418573670544, 12.76057685722273, 11.81774503966809, 16.538260870614, 22.20221269367999),
        @(15.62745848615193, 6.238654370285997, 7.476443217733217, 0, 40.8793177260399, 48.24658325200691, 19.10585748447998, 30.11644858810017, 10.45029001747888, 12.79797466623664, 11.81697199334803, 16.54692589639859, 22.18646506259116),
        @(15.61639012757499, 6.23184291105869, 7.475738948493041, 0, 40.88060018877849, 48.24656291286845, 19.1077892545266, 30.1197861614882, 10.45112539489982, 12.78991503812417, 11.81712964194006, 16.54504529268785, 22.18984478677815),
        @(15.58030833985485, 6.209554358940252, 7.473461554905018, 0, 40.8849752272894, 48.24676937526432, 19.11415725124809, 30.13078866755088, 10.45386359484109, 12.76365114612853, 11.81767745556028, 16.53896723566143, 22.20091180321585),
        @(15.55821917665542, 6.195844670280326, 7.472081510060212, 0, 40.88780311260677, 48.24710547904194, 19.11811014371983, 30.13761865046756, 10.45555147225047, 12.74757991561098, 11.81803890920864, 16.5352866002758, 22.20772526983988),
        @(15.43185196341527, 6.11641217054805, 7.464403064789306, 0, 40.90627820077419, 48.25224541352318, 19.14156437688262, 30.17814731195003, 10.4653875242327, 12.65575294870936, 11.82050791344332, 16.51484700702285, 22.24730162368662),
        @(15.35455966291009, 6.066920329995503, 7.459899223613641, 0, 40.91963764882298, 48.2582668019202, 19.15666827359954, 30.20424920840436, 10.47156780897769, 12.59968676940874, 11.82237810128902, 16.50289410316168, 22.2720552623075),
        @(15.31020214219872, 6.038178503815729, 7.457385206918524, 0, 40.92806612983611, 48.26278424556657, 19.16561857807694, 30.21971745937894, 10.47517638070752, 12.56754694710155, 11.8235847107524, 16.49623596741739, 22.28646775194724),
        @(15.29520182302323, 6.028400180057303, 7.456547188423595, 0, 40.93104773692451, 48.26449481009882, 19.16869415818084, 30.22503289754266, 10.47640743810381, 12.55668444193464, 11.82401575979054, 16.49401901813651, 22.29137763792313),
        @(15.3627776448874, 6.072217437590494, 7.460370762396876, 0, 40.9181384629323, 48.25751671842257, 19.15503322447821, 30.20142350912692, 10.47090433790229, 12.60564416077288, 11.82216547213796, 16.50414409853854, 22.26940210674485),
        @(15.59090401980434, 6.216113014495327, 7.474127373454134, 0, 40.88365932062227, 48.24666504549847, 19.11227591483338, 30.12753807678924, 10.45305710013005, 12.77136217487324, 11.81751115305079, 16.54074369117964, 22.19765394365903),
        @(15.74057261900824, 6.307604889364218, 7.483787554218423, 0, 40.86775180217187, 48.24895614027337, 19.0866728573091, 30.08330567198161, 10.44186990940232, 12.88041863196817, 11.8156293328178, 16.56656310926908, 22.15231191318492),
        @(15.66065379844078, 6.25901207488555, 7.478571053185538, 0, 40.87563584052995, 48.24687512801881, 19.10012347307593, 30.10654207398271, 10.44779719200017, 12.82215470722855, 11.81652784389007, 16.55261058010774, 22.17637038513248),
        @(15.35906204617795, 6.069823512163437, 7.46015734504453, 0, 40.91881391499671, 48.25785254279437, 19.15577159887355, 30.202699568922, 10.47120412044733, 12.60295053092947, 11.82226119237695, 16.50357830947894, 22.27060103355327),
        @(15.03771393662499, 5.855292681511812, 7.443224925723217, 0, 40.9938484350244, 48.30998544719665, 19.22584981914112, 30.32381929840816, 10.49851202268439, 12.37075505089441, 11.83345404899817, 16.45899583788475, 22.37893378638985)
)

$numCols = 13  # columns B..N
$data = New-Object "object[,]" $newRows.Count, $numCols
for ($r = 0; $r -lt $newRows.Count; $r++) {
    for ($c = 0; $c -lt $numCols; $c++) {
        $data[$r, $c] = $newRows[$r][$c]
    }
}

$ws.Range("B2:N25").Value = $data
